# TC-73 added in requirement tab
# Adds a new worksheet "tc073" (positioned right after "tc064", right
# before "tc047") holding the requirement data for TC073, mirroring the
# layout used by the neighbouring tc063/tc064/tc047 requirement sheets.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately before "tc047" -- i.e. right after
# "tc064" -- which matches the sheet order in the target workbook.
$beforeSheet = $wb.Worksheets.Item("tc047")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "tc073"

# Header row
$ws.Range("A1").Value = "ModuleName"
$ws.Range("B1").Value = "Requirementitle "
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Priority"
$ws.Range("E1").Value = "Status"
$ws.Range("F1").Value = "Type"

# Data row for TC073
$ws.Range("A2").Value = "Epic Mohit"
$ws.Range("B2").Value = "Notification for New Requirement 8-01-2026"
$ws.Range("C2").Value = "delete notification Flow testing"
$ws.Range("D2").Value = "Must have"
$ws.Range("E2").Value = "New"
$ws.Range("F2").Value = "Functional"

# Match the column widths used on the sibling requirement sheets and
# select the populated range, mirroring the new sheet's saved view state.
$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(2).ColumnWidth = 51.67
$ws.Columns.Item(3).ColumnWidth = 26.67

$ws.Range("A1:F2").Select()
